$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 74, shifting rows 74:94 down to 75:95.
$ws.Rows.Item(74).Insert()

# Populate the new row 74 with the new market record.
$ws.Range("A74").Value = 10
$ws.Range("B74").Value = "Vega Modelo de Temuco"
$ws.Range("C74").Value = "La Araucanía"
$ws.Range("D74").Value = 44754
$ws.Range("D74").NumberFormat = $ws.Range("D75").NumberFormat
$ws.Range("E74").Value = 9
$ws.Range("F74").Value = 100112035
$ws.Range("G74").Value = "Bruselas (repollito)"
$ws.Range("H74").Value = "Sin especificar"
$ws.Range("I74").Value = "Primera"
$ws.Range("J74").Value = 40
$ws.Range("K74").Value = 26000
$ws.Range("L74").Value = 26000
$ws.Range("M74").Value = 26000
$ws.Range("N74").Value = "$/malla 10 kilos"
$ws.Range("O74").Value = "Región Metropolitana"
$ws.Range("P74").Value = 2600
$ws.Range("Q74").Value = 10
$ws.Range("R74").Value = "Hortaliza"
